# "Actualizada la planilla de sistemas" - se agregan las tarjetas magneticas
# cargadas recientemente en el sistema (hoja "Tarjetas").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tarjetas")
$ws.Activate()

# Fila 24: el portador de la tarjeta ya no se identifica como
# "Bombero 01 (Fernando Sánchez)" sino simplemente por su nombre, y pasa a
# pertenecer al área "Espacios Físicos".
$ws.Cells.Item(24, 1).Value = "Fernando Sánchez"
$ws.Cells.Item(24, 2).Value = 32554035
$ws.Cells.Item(24, 3).Value = "Espacios Físicos"
$ws.Cells.Item(24, 4).Value = "0008382231 | 127,59159"

# Nuevas tarjetas magnéticas cargadas recientemente (filas 25-28).
$ws.Cells.Item(25, 1).Value = "Andrés Habermehl"
$ws.Cells.Item(25, 2).Value = 35970957
$ws.Cells.Item(25, 3).Value = "Espacios Físicos"
$ws.Cells.Item(25, 4).Value = "0008380410 | 127,57338"

$ws.Cells.Item(26, 1).Value = "Ezequiel Mahafud"
$ws.Cells.Item(26, 2).Value = 38999220
$ws.Cells.Item(26, 3).Value = "Espacios Físicos"
$ws.Cells.Item(26, 4).Value = "0008385193 | 127,62121"

$ws.Cells.Item(27, 1).Value = "Demian Szmulewicz"
$ws.Cells.Item(27, 2).Value = 26352654
$ws.Cells.Item(27, 3).Value = "Espacios Físicos"
$ws.Cells.Item(27, 4).Value = "0008380715 | 127,57643"

$ws.Cells.Item(28, 1).Value = "Jesus Juarez"
$ws.Cells.Item(28, 2).Value = 33194188
$ws.Cells.Item(28, 3).Value = "Espacios Físicos"
$ws.Cells.Item(28, 4).Value = "0008383067 | 127,59995"

# Vista: desplazar la ventana y dejar el cursor listo para la próxima fila.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D29").Select()

# Configuración de impresión tocada al revisar la hoja.
$ws.PageSetup.Orientation = 1
